# Add two new rows (9 and 10) of localization data to the sheet, mirroring
# the existing table of "filename / line number / english / translated /
# converted" script-line entries.
#
# Values are entered column-by-column (C, then A, then D, then E, then B)
# to reproduce the original shared-string insertion order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C (Англ. Строка / English line) ---
# C9 repeats the English line already used in rows 2 and 6.
$ws.Range("C9").Value = " Oh, [hero] and\n[partner]!"
$ws.Range("C10").Value = ' We\''re counting on you![K]\nGood luck!'

# --- Column A (Название файла в скриптах и цвет / script filename) ---
$ws.Range("A9").Value = "SCRIPT/T01P02A/um1408.ssb"

# --- Column D (Переведенная строка / translated Russian line) ---
$ws.Range("D9").Value = " Ой, [hero] и [partner]!"
$ws.Range("D10").Value = " Мы надеемся на вас![K] Удачи!"

# --- Column E (Конвертированная строка / re-encoded/garbled line) ---
$ws.Range("E9").Value = " Ïê, [hero] é [partner]!"
$ws.Range("E10").Value = " Íú îàäååíòÿ îà âàò![K] Ôäàœé!"

# --- Column B (Номер строки / line number) ---
$ws.Range("B9").Value = 218
$ws.Range("B10").Value = 221

# Row heights to match the wrapped, two-line content (as in rows 2/6/8 and 7).
$ws.Rows.Item(9).RowHeight = 43.2
$ws.Rows.Item(10).RowHeight = 26.4

# Leave the selection on the newly entered translation cell, as in the
# edited workbook.
$ws.Range("D9").Select()
